# Update the nowcast table with the 2025Q4 vintage data.
# Only the top block of data rows (2-7) is refreshed; the header (row 1)
# and the trailing historical rows (8-11) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-09-30", 0.19453243892083816, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("2025-10-15", 0.25614714246057269, 0, 0.00072863431194951963, 0.00099207090997436415, -0.0001343652788438679, 0.00063524561232872689, -0.00007480431668907172, -0.00046280605566136355, 0, -0.00047170149079100998),
    @("2025-10-30", 0.3110260838382391, 0.040070426296224218, 0, -0.000083617729564361298, -0.00001490992572603138, 0, 0.00012064226663808732, -0.0087502047124025633, 0.0025363668788354332, -0.000060375720833838464),
    @("2025-11-15", 0.30794307504700436, 0, -0.00064889452771484471, 0.000016665468506536254, 0.0071358579926211335, -0.0015648006627433163, 0.00066085543953165371, -0.00079303136428770951, 0, -0.011257684745037855),
    @("2025-11-30", 0.3316983548611655, -0.01773461821504857, 0, -0.00021478259595248021, 0.00042812033183026968, 0, 0.00056165252309212511, -0.0062917667795795214, 0, 0.00018810843477501393),
    @("2025-12-15", 0.24825101149903989, 0, -0.0184257860329703, -0.0030530078896552231, 0.0071358999548796981, 0.0028458481366997338, 0, 0, 0, 0.01048649352055564)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]

    # Column A holds a date label that must stay plain text (matching the
    # workbook's existing convention), so force text formatting before
    # assigning it and restore the default style afterward.
    $dateCell = $ws.Cells.Item($rowIndex, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rowValues[0]
    $dateCell.Style = "Normal"

    for ($col = 2; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
